# Update master to output generated at 4250d90
$d = $word.ActiveDocument

# 1. Update the date heading at the top of the document.
$d.Content.Find.Execute("2024-05-28 Tuesday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2024-05-29 Wednesday", 2)

# 2. Update each arithmetic-fact cell in the (single) table.
#    Cells are addressed directly by (row, column) so that values which are
#    reused as both an old and a new answer elsewhere in the table never
#    collide with a text-based Find/Replace.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "94÷6=15, 4"
$t.Cell(1, 2).Range.Text = "57÷2=28, 1"
$t.Cell(1, 3).Range.Text = "32÷7=4, 4"
$t.Cell(1, 4).Range.Text = "33÷8=4, 1"
$t.Cell(1, 5).Range.Text = "46÷7=6, 4"

$t.Cell(5, 1).Range.Text = "91÷5=18, 1"
$t.Cell(5, 2).Range.Text = "88÷3=29, 1"
$t.Cell(5, 3).Range.Text = "63÷5=12, 3"
$t.Cell(5, 4).Range.Text = "50÷4=12, 2"
$t.Cell(5, 5).Range.Text = "51÷2=25, 1"

$t.Cell(9, 1).Range.Text = "61÷2=30, 1"
$t.Cell(9, 2).Range.Text = "21÷2=10, 1"
$t.Cell(9, 3).Range.Text = "93÷8=11, 5"
$t.Cell(9, 4).Range.Text = "70÷4=17, 2"
$t.Cell(9, 5).Range.Text = "86÷5=17, 1"

$t.Cell(13, 1).Range.Text = "46÷3=15, 1"
$t.Cell(13, 2).Range.Text = "62÷4=15, 2"
$t.Cell(13, 3).Range.Text = "54÷8=6, 6"
$t.Cell(13, 4).Range.Text = "84÷5=16, 4"
$t.Cell(13, 5).Range.Text = "17÷2=8, 1"

$t.Cell(17, 1).Range.Text = "60÷6=10, 0"
$t.Cell(17, 2).Range.Text = "66÷6=11, 0"
$t.Cell(17, 3).Range.Text = "15÷7=2, 1"
$t.Cell(17, 4).Range.Text = "19÷3=6, 1"
$t.Cell(17, 5).Range.Text = "33÷7=4, 5"
